$d = $word.ActiveDocument

$newValues = @(
    "28-11=17",
    "6+8=14",
    "42+21=63",
    "80-34=46",
    "7+76=83",
    "36+4=40",
    "50+16=66",
    "78-3=75",
    "63+17=80",
    "89+5=94",
    "40+7=47",
    "28-16=12",
    "8+45=53",
    "97-91=6",
    "40+4=44",
    "81-68=13",
    "87+9=96",
    "98-82=16",
    "55-40=15",
    "59-41=18",
    "98-60=38",
    "92-34=58",
    "37-19=18",
    "56-46=10",
    "87-20=67",
    "48-3=45",
    "94-9=85",
    "40+0=40",
    "17-2=15",
    "58-52=6",
    "36-33=3",
    "78-62=16",
    "87-15=72",
    "30+69=99",
    "9-0=9",
    "93-35=58",
    "45+32=77",
    "85-57=28",
    "54-20=34",
    "41+47=88",
    "62-4=58",
    "69-28=41",
    "85-5=80",
    "67-39=28",
    "41+57=98",
    "55+17=72",
    "7+69=76",
    "99-67=32",
    "90-46=44",
    "17+6=23",
    "90-14=76",
    "41+9=50",
    "77+15=92",
    "20+75=95",
    "62-14=48",
    "41+33=74",
    "70-58=12",
    "91-33=58",
    "21-9=12",
    "7+65=72",
    "82-74=8",
    "52+15=67",
    "67-7=60",
    "75+12=87",
    "17+61=78",
    "20+73=93",
    "26-5=21",
    "2+76=78",
    "29+22=51",
    "95-92=3",
    "75+17=92",
    "70-0=70",
    "92-41=51",
    "39+47=86",
    "95-62=33",
    "95-44=51",
    "55-36=19",
    "88-82=6",
    "71+22=93",
    "63-28=35",
    "97-91=6",
    "57+38=95",
    "5+76=81",
    "81-57=24",
    "12-11=1",
    "66-35=31",
    "23+37=60",
    "64-14=50",
    "17+64=81",
    "33+55=88",
    "90-57=33",
    "52+28=80",
    "69+25=94",
    "77+12=89",
    "90-66=24",
    "97-75=22",
    "49+18=67",
    "52-5=47",
    "49-26=23",
    "41-3=38"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        if ($idx -lt $newValues.Count) {
            $cell = $t.Cell($r, $c)
            $cell.Range.Text = $newValues[$idx]
        }
        $idx = $idx + 1
    }
}

Write-Output "Updated $idx cells (rows=$rows cols=$cols)"
